$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = -0.21823451781839509
$ws.Cells.Item(1, 2).Value = 0.21804058587594
$ws.Cells.Item(2, 1).Value = -0.15727086447088112
$ws.Cells.Item(2, 2).Value = 0.15685396304133015
$ws.Cells.Item(3, 1).Value = -0.1071485688817031
$ws.Cells.Item(3, 2).Value = 0.10697114947461372
$ws.Cells.Item(4, 1).Value = -0.09897114953202113
$ws.Cells.Item(4, 2).Value = 0.098581280913808911
$ws.Cells.Item(5, 1).Value = -0.095581280946419156
$ws.Cells.Item(5, 2).Value = 0.094262424806784928
$ws.Cells.Item(6, 1).Value = -0.040336610583187493
$ws.Cells.Item(6, 2).Value = 0.04001328858917752
$ws.Cells.Item(7, 1).Value = -0.030013288672153582
$ws.Cells.Item(7, 2).Value = 0.029944494975937541
$ws.Cells.Item(8, 1).Value = -0.019944495061059886
$ws.Cells.Item(8, 2).Value = 0.019848439018641084
$ws.Cells.Item(9, 1).Value = -0.0178484390612903
$ws.Cells.Item(9, 2).Value = 0.017777640625388003
$ws.Cells.Item(10, 1).Value = -0.015777640670300741
$ws.Cells.Item(10, 2).Value = 0.01577469520703545
$ws.Cells.Item(11, 1).Value = -0.012774695257741442
$ws.Cells.Item(11, 2).Value = 0.012769695076466192
$ws.Cells.Item(12, 1).Value = -0.0092696951303925523
$ws.Cells.Item(12, 2).Value = 0.0092408823871736878
$ws.Cells.Item(13, 1).Value = -0.0057408824429137084
$ws.Cells.Item(13, 2).Value = 0.0057339858633591589
$ws.Cells.Item(14, 1).Value = 0.0022660140551131747
$ws.Cells.Item(14, 2).Value = -0.002266056368161351
$ws.Cells.Item(15, 1).Value = 0.0032660563259723219
$ws.Cells.Item(15, 2).Value = -0.0032679850042036662
$ws.Cells.Item(16, 1).Value = -0.0060339638241582172
$ws.Cells.Item(16, 2).Value = 0.006003423258485352
$ws.Cells.Item(17, 1).Value = -0.0040034233071857273
$ws.Cells.Item(17, 2).Value = 0.0039999999397943853
$ws.Cells.Item(18, 1).Value = -0.054769721440276697
$ws.Cells.Item(18, 2).Value = 0.054666415890579856
$ws.Cells.Item(19, 1).Value = -0.050666415915244567
$ws.Cells.Item(19, 2).Value = 0.049910654173991986
$ws.Cells.Item(20, 1).Value = -0.045910654207485635
$ws.Cells.Item(20, 2).Value = 0.045697554753706626
$ws.Cells.Item(21, 1).Value = -0.0040057988573156678
$ws.Cells.Item(21, 2).Value = 0.0039999999645248252
$ws.Cells.Item(22, 1).Value = -0.045705394191250193
$ws.Cells.Item(22, 2).Value = 0.045494112976712131
$ws.Cells.Item(23, 1).Value = -0.040494113015801858
$ws.Cells.Item(23, 2).Value = 0.040098037145776644
$ws.Cells.Item(24, 1).Value = -0.020098037272974878
$ws.Cells.Item(24, 2).Value = 0.019999999871133767
$ws.Cells.Item(25, 1).Value = -0.051925814260393111
$ws.Cells.Item(25, 2).Value = 0.051889201174555311
$ws.Cells.Item(26, 1).Value = -0.049389201214236067
$ws.Cells.Item(26, 2).Value = 0.049345350170680646
$ws.Cells.Item(27, 1).Value = -0.046845350211541792
$ws.Cells.Item(27, 2).Value = 0.04660151496849041
$ws.Cells.Item(28, 1).Value = -0.044601515011974513
$ws.Cells.Item(28, 2).Value = 0.044449010147412871
$ws.Cells.Item(29, 1).Value = -0.037449010222869283
$ws.Cells.Item(29, 2).Value = 0.037415607444127197
$ws.Cells.Item(30, 1).Value = 0.022584392180328816
$ws.Cells.Item(30, 2).Value = -0.022618109256106411
$ws.Cells.Item(31, 1).Value = 0.029618109180864494
$ws.Cells.Item(31, 2).Value = -0.029641039974835692
$ws.Cells.Item(32, 1).Value = -0.004001239050550609
$ws.Cells.Item(32, 2).Value = 0.0039999999423798727
